$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98, pushing the existing rows 98-106 down
# to 99-107 (this also grows the sheet's used range to A1:T107).
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly price record.
$ws.Range("A98").Value = 2
$ws.Range("B98").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C98").Value = "Coquimbo"
$ws.Range("D98").Value = 44776
$ws.Range("E98").Value = 4
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100107
$ws.Range("H98").Value = "Otros"
$ws.Range("I98").Value = 100107011
$ws.Range("J98").Value = "Tuna"
$ws.Range("K98").Value = "Sin especificar"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 200
$ws.Range("N98").Value = 20000
$ws.Range("O98").Value = 21000
$ws.Range("P98").Value = 20500
$ws.Range("Q98").Value = "$/caja 18 kilos"
$ws.Range("R98").Value = "Provincia de Limarí"
$ws.Range("S98").Value = 1139
$ws.Range("T98").Value = 18
